$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = New-Object 'object[,]' 36,16
$data[0,0] = 3
$data[0,1] = 1
$data[0,2] = 1.141514
$data[0,3] = 3.424542
$data[0,4] = 0.1286876091642341
$data[0,5] = 0.1286876091642341
$data[0,6] = 3
$data[0,7] = 1
$data[0,8] = 3.264166666666667
$data[0,9] = 9.7925
$data[0,10] = 0.02701285365228691
$data[0,11] = 0.02701285365228691
$data[0,12] = 3.726091948333334
$data[0,13] = 33.53482753500001
$data[0,14] = 0.003476219553216153
$data[0,15] = 0.003476219553216152
$data[1,0] = 3
$data[1,1] = 1
$data[1,2] = 1.141514
$data[1,3] = 3.424542
$data[1,4] = 0.1286876091642341
$data[1,5] = 0.1286876091642341
$data[1,6] = 3
$data[1,7] = 1
$data[1,8] = 37.12743
$data[1,9] = 111.38229
$data[1,10] = 0.3072508041078968
$data[1,11] = 0.3072508041078968
$data[1,12] = 42.38148112902
$data[1,13] = 381.43333016118
$data[1,14] = 0.03953937139443369
$data[1,15] = 0.03953937139443368
$data[2,0] = 3
$data[2,1] = 1
$data[2,2] = 1.141514
$data[2,3] = 3.424542
$data[2,4] = 0.1286876091642341
$data[2,5] = 0.1286876091642341
$data[2,6] = 3
$data[2,7] = 1
$data[2,8] = 15.075014
$data[2,9] = 45.225042
$data[2,10] = 0.1247543978518794
$data[2,11] = 0.1247543978518794
$data[2,12] = 17.208339531196
$data[2,13] = 154.875055780764
$data[2,14] = 0.01605434519228203
$data[2,15] = 0.01605434519228202
$data[3,0] = 3
$data[3,1] = 1
$data[3,2] = 1.141514
$data[3,3] = 3.424542
$data[3,4] = 0.1286876091642341
$data[3,5] = 0.1286876091642341
$data[3,6] = 3
$data[3,7] = 1
$data[3,8] = 31.08386466666667
$data[3,9] = 93.25159400000001
$data[3,10] = 0.2572368303869774
$data[3,11] = 0.2572368303869774
$data[3,12] = 35.48266669110534
$data[3,13] = 319.344000219948
$data[3,14] = 0.03310319269148574
$data[3,15] = 0.03310319269148573
$data[4,0] = 3
$data[4,1] = 1
$data[4,2] = 1.141514
$data[4,3] = 3.424542
$data[4,4] = 0.1286876091642341
$data[4,5] = 0.1286876091642341
$data[4,6] = 3
$data[4,7] = 1
$data[4,8] = 5.556591666666667
$data[4,9] = 16.669775
$data[4,10] = 0.04598398697896871
$data[4,11] = 0.04598398697896871
$data[4,12] = 6.342927179783334
$data[4,13] = 57.08634461805001
$data[4,14] = 0.005917569344162757
$data[4,15] = 0.005917569344162756
$data[5,0] = 3
$data[5,1] = 1
$data[5,2] = 1.141514
$data[5,3] = 3.424542
$data[5,4] = 0.1286876091642341
$data[5,5] = 0.1286876091642341
$data[5,6] = 3
$data[5,7] = 1
$data[5,8] = 28.73046866666667
$data[5,9] = 86.191406
$data[5,10] = 0.2377611270219907
$data[5,11] = 0.2377611270219907
$data[5,12] = 32.79623220956134
$data[5,13] = 295.166089886052
$data[5,14] = 0.03059691098865377
$data[5,15] = 0.03059691098865376
$data[6,0] = 3
$data[6,1] = 1
$data[6,2] = 1.348859666666667
$data[6,3] = 4.046578999999999
$data[6,4] = 0.1520625464088912
$data[6,5] = 0.1520625464088912
$data[6,6] = 3
$data[6,7] = 1
$data[6,8] = 3.264166666666667
$data[6,9] = 9.7925
$data[6,10] = 0.02701285365228691
$data[6,11] = 0.02701285365228691
$data[6,12] = 4.402902761944445
$data[6,13] = 39.6261248575
$data[6,14] = 0.004107643312137466
$data[6,15] = 0.004107643312137466
$data[7,0] = 3
$data[7,1] = 1
$data[7,2] = 1.348859666666667
$data[7,3] = 4.046578999999999
$data[7,4] = 0.1520625464088912
$data[7,5] = 0.1520625464088912
$data[7,6] = 3
$data[7,7] = 1
$data[7,8] = 37.12743
$data[7,9] = 111.38229
$data[7,10] = 0.3072508041078968
$data[7,11] = 0.3072508041078968
$data[7,12] = 50.07969285398999
$data[7,13] = 450.7172356859099
$data[7,14] = 0.04672133965882622
$data[7,15] = 0.04672133965882622
$data[8,0] = 3
$data[8,1] = 1
$data[8,2] = 1.348859666666667
$data[8,3] = 4.046578999999999
$data[8,4] = 0.1520625464088912
$data[8,5] = 0.1520625464088912
$data[8,6] = 3
$data[8,7] = 1
$data[8,8] = 15.075014
$data[8,9] = 45.225042
$data[8,10] = 0.1247543978518794
$data[8,11] = 0.1247543978518794
$data[8,12] = 20.33407835903533
$data[8,13] = 183.006705231318
$data[8,14] = 0.01897047141306469
$data[8,15] = 0.01897047141306469
$data[9,0] = 3
$data[9,1] = 1
$data[9,2] = 1.348859666666667
$data[9,3] = 4.046578999999999
$data[9,4] = 0.1520625464088912
$data[9,5] = 0.1520625464088912
$data[9,6] = 3
$data[9,7] = 1
$data[9,8] = 31.08386466666667
$data[9,9] = 93.25159400000001
$data[9,10] = 0.2572368303869774
$data[9,11] = 0.2572368303869774
$data[9,12] = 41.92777133299178
$data[9,13] = 377.349941996926
$data[9,14] = 0.03911608745879584
$data[9,15] = 0.03911608745879584
$data[10,0] = 3
$data[10,1] = 1
$data[10,2] = 1.348859666666667
$data[10,3] = 4.046578999999999
$data[10,4] = 0.1520625464088912
$data[10,5] = 0.1520625464088912
$data[10,6] = 3
$data[10,7] = 1
$data[10,8] = 5.556591666666667
$data[10,9] = 16.669775
$data[10,10] = 0.04598398697896871
$data[10,11] = 0.04598398697896871
$data[10,12] = 7.495062383302778
$data[10,13] = 67.455561449725
$data[10,14] = 0.00699244215405528
$data[10,15] = 0.00699244215405528
$data[11,0] = 3
$data[11,1] = 1
$data[11,2] = 1.348859666666667
$data[11,3] = 4.046578999999999
$data[11,4] = 0.1520625464088912
$data[11,5] = 0.1520625464088912
$data[11,6] = 3
$data[11,7] = 1
$data[11,8] = 28.73046866666667
$data[11,9] = 86.191406
$data[11,10] = 0.2377611270219907
$data[11,11] = 0.2377611270219907
$data[11,12] = 38.75337038889711
$data[11,13] = 348.7803335000739
$data[11,14] = 0.03615456241201175
$data[11,15] = 0.03615456241201175
$data[12,0] = 3
$data[12,1] = 1
$data[12,2] = 2.436302666666667
$data[12,3] = 7.308908000000001
$data[12,4] = 0.2746545074118945
$data[12,5] = 0.2746545074118945
$data[12,6] = 3
$data[12,7] = 1
$data[12,8] = 3.264166666666667
$data[12,9] = 9.7925
$data[12,10] = 0.02701285365228691
$data[12,11] = 0.02701285365228691
$data[12,12] = 7.952497954444445
$data[12,13] = 71.57248159000001
$data[12,14] = 0.007419202013658457
$data[12,15] = 0.007419202013658457
$data[13,0] = 3
$data[13,1] = 1
$data[13,2] = 2.436302666666667
$data[13,3] = 7.308908000000001
$data[13,4] = 0.2746545074118945
$data[13,5] = 0.2746545074118945
$data[13,6] = 3
$data[13,7] = 1
$data[13,8] = 37.12743
$data[13,9] = 111.38229
$data[13,10] = 0.3072508041078968
$data[13,11] = 0.3072508041078968
$data[13,12] = 90.45365671548
$data[13,13] = 814.08291043932
$data[13,14] = 0.0843878182541629
$data[13,15] = 0.0843878182541629
$data[14,0] = 3
$data[14,1] = 1
$data[14,2] = 2.436302666666667
$data[14,3] = 7.308908000000001
$data[14,4] = 0.2746545074118945
$data[14,5] = 0.2746545074118945
$data[14,6] = 3
$data[14,7] = 1
$data[14,8] = 15.075014
$data[14,9] = 45.225042
$data[14,10] = 0.1247543978518794
$data[14,11] = 0.1247543978518794
$data[14,12] = 36.72729680823733
$data[14,13] = 330.5456712741361
$data[14,14] = 0.03426435768947545
$data[14,15] = 0.03426435768947545
$data[15,0] = 3
$data[15,1] = 1
$data[15,2] = 2.436302666666667
$data[15,3] = 7.308908000000001
$data[15,4] = 0.2746545074118945
$data[15,5] = 0.2746545074118945
$data[15,6] = 3
$data[15,7] = 1
$data[15,8] = 31.08386466666667
$data[15,9] = 93.25159400000001
$data[15,10] = 0.2572368303869774
$data[15,11] = 0.2572368303869774
$data[15,12] = 75.72970237770579
$data[15,13] = 681.5673213993522
$data[15,14] = 0.07065125493813235
$data[15,15] = 0.07065125493813235
$data[16,0] = 3
$data[16,1] = 1
$data[16,2] = 2.436302666666667
$data[16,3] = 7.308908000000001
$data[16,4] = 0.2746545074118945
$data[16,5] = 0.2746545074118945
$data[16,6] = 3
$data[16,7] = 1
$data[16,8] = 5.556591666666667
$data[16,9] = 16.669775
$data[16,10] = 0.04598398697896871
$data[16,11] = 0.04598398697896871
$data[16,12] = 13.53753909507778
$data[16,13] = 121.8378518557
$data[16,14] = 0.01262970929254362
$data[16,15] = 0.01262970929254362
$data[17,0] = 3
$data[17,1] = 1
$data[17,2] = 2.436302666666667
$data[17,3] = 7.308908000000001
$data[17,4] = 0.2746545074118945
$data[17,5] = 0.2746545074118945
$data[17,6] = 3
$data[17,7] = 1
$data[17,8] = 28.73046866666667
$data[17,9] = 86.191406
$data[17,10] = 0.2377611270219907
$data[17,11] = 0.2377611270219907
$data[17,12] = 69.99611742718311
$data[17,13] = 629.9650568446481
$data[17,14] = 0.06530216522392174
$data[17,15] = 0.06530216522392174
$data[18,0] = 1
$data[18,1] = 0.3333333333333333
$data[18,2] = 0.4110173333333333
$data[18,3] = 1.233052
$data[18,4] = 0.04633568922652349
$data[18,5] = 0.04633568922652349
$data[18,6] = 3
$data[18,7] = 1
$data[18,8] = 3.264166666666667
$data[18,9] = 9.7925
$data[18,10] = 0.02701285365228691
$data[18,11] = 0.02701285365228691
$data[18,12] = 1.341629078888889
$data[18,13] = 12.07466171
$data[18,14] = 0.001251659191953926
$data[18,15] = 0.001251659191953926
$data[19,0] = 1
$data[19,1] = 0.3333333333333333
$data[19,2] = 0.4110173333333333
$data[19,3] = 1.233052
$data[19,4] = 0.04633568922652349
$data[19,5] = 0.04633568922652349
$data[19,6] = 3
$data[19,7] = 1
$data[19,8] = 37.12743
$data[19,9] = 111.38229
$data[19,10] = 0.3072508041078968
$data[19,11] = 0.3072508041078968
$data[19,12] = 15.26001727212
$data[19,13] = 137.34015544908
$data[19,14] = 0.01423667777374295
$data[19,15] = 0.01423667777374295
$data[20,0] = 1
$data[20,1] = 0.3333333333333333
$data[20,2] = 0.4110173333333333
$data[20,3] = 1.233052
$data[20,4] = 0.04633568922652349
$data[20,5] = 0.04633568922652349
$data[20,6] = 3
$data[20,7] = 1
$data[20,8] = 15.075014
$data[20,9] = 45.225042
$data[20,10] = 0.1247543978518794
$data[20,11] = 0.1247543978518794
$data[20,12] = 6.196092054242667
$data[20,13] = 55.764828488184
$data[20,14] = 0.005780581008506753
$data[20,15] = 0.005780581008506753
$data[21,0] = 1
$data[21,1] = 0.3333333333333333
$data[21,2] = 0.4110173333333333
$data[21,3] = 1.233052
$data[21,4] = 0.04633568922652349
$data[21,5] = 0.04633568922652349
$data[21,6] = 3
$data[21,7] = 1
$data[21,8] = 31.08386466666667
$data[21,9] = 93.25159400000001
$data[21,10] = 0.2572368303869774
$data[21,11] = 0.2572368303869774
$data[21,12] = 12.77600716498756
$data[21,13] = 114.984064484888
$data[21,14] = 0.01191924583042692
$data[21,15] = 0.01191924583042692
$data[22,0] = 1
$data[22,1] = 0.3333333333333333
$data[22,2] = 0.4110173333333333
$data[22,3] = 1.233052
$data[22,4] = 0.04633568922652349
$data[22,5] = 0.04633568922652349
$data[22,6] = 3
$data[22,7] = 1
$data[22,8] = 5.556591666666667
$data[22,9] = 16.669775
$data[22,10] = 0.04598398697896871
$data[22,11] = 0.04598398697896871
$data[22,12] = 2.283855489255556
$data[22,13] = 20.5546994033
$data[22,14] = 0.002130699730053997
$data[22,15] = 0.002130699730053997
$data[23,0] = 1
$data[23,1] = 0.3333333333333333
$data[23,2] = 0.4110173333333333
$data[23,3] = 1.233052
$data[23,4] = 0.04633568922652349
$data[23,5] = 0.04633568922652349
$data[23,6] = 3
$data[23,7] = 1
$data[23,8] = 28.73046866666667
$data[23,9] = 86.191406
$data[23,10] = 0.2377611270219907
$data[23,11] = 0.2377611270219907
$data[23,12] = 11.80872061679022
$data[23,13] = 106.278485551112
$data[23,14] = 0.01101682569183894
$data[23,15] = 0.01101682569183894
$data[24,0] = 3
$data[24,1] = 1
$data[24,2] = 2.206596333333333
$data[24,3] = 6.619789
$data[24,4] = 0.2487587594433638
$data[24,5] = 0.2487587594433639
$data[24,6] = 3
$data[24,7] = 1
$data[24,8] = 3.264166666666667
$data[24,9] = 9.7925
$data[24,10] = 0.02701285365228691
$data[24,11] = 0.02701285365228691
$data[24,12] = 7.202698198055555
$data[24,13] = 64.8242837825
$data[24,14] = 0.006719683963568032
$data[24,15] = 0.006719683963568033
$data[25,0] = 3
$data[25,1] = 1
$data[25,2] = 2.206596333333333
$data[25,3] = 6.619789
$data[25,4] = 0.2487587594433638
$data[25,5] = 0.2487587594433639
$data[25,6] = 3
$data[25,7] = 1
$data[25,8] = 37.12743
$data[25,9] = 111.38229
$data[25,10] = 0.3072508041078968
$data[25,11] = 0.3072508041078968
$data[25,12] = 81.92525090408999
$data[25,13] = 737.32725813681
$data[25,14] = 0.07643132886785642
$data[25,15] = 0.07643132886785643
$data[26,0] = 3
$data[26,1] = 1
$data[26,2] = 2.206596333333333
$data[26,3] = 6.619789
$data[26,4] = 0.2487587594433638
$data[26,5] = 0.2487587594433639
$data[26,6] = 3
$data[26,7] = 1
$data[26,8] = 15.075014
$data[26,9] = 45.225042
$data[26,10] = 0.1247543978518794
$data[26,11] = 0.1247543978518794
$data[26,12] = 33.26447061734866
$data[26,13] = 299.380235556138
$data[26,14] = 0.03103374924473738
$data[26,15] = 0.03103374924473738
$data[27,0] = 3
$data[27,1] = 1
$data[27,2] = 2.206596333333333
$data[27,3] = 6.619789
$data[27,4] = 0.2487587594433638
$data[27,5] = 0.2487587594433639
$data[27,6] = 3
$data[27,7] = 1
$data[27,8] = 31.08386466666667
$data[27,9] = 93.25159400000001
$data[27,10] = 0.2572368303869774
$data[27,11] = 0.2572368303869774
$data[27,12] = 68.58954179929623
$data[27,13] = 617.3058761936661
$data[27,14] = 0.06398991481020751
$data[27,15] = 0.06398991481020751
$data[28,0] = 3
$data[28,1] = 1
$data[28,2] = 2.206596333333333
$data[28,3] = 6.619789
$data[28,4] = 0.2487587594433638
$data[28,5] = 0.2487587594433639
$data[28,6] = 3
$data[28,7] = 1
$data[28,8] = 5.556591666666667
$data[28,9] = 16.669775
$data[28,10] = 0.04598398697896871
$data[28,11] = 0.04598398697896871
$data[28,12] = 12.26115479749722
$data[28,13] = 110.350393177475
$data[28,14] = 0.01143891955514805
$data[28,15] = 0.01143891955514805
$data[29,0] = 3
$data[29,1] = 1
$data[29,2] = 2.206596333333333
$data[29,3] = 6.619789
$data[29,4] = 0.2487587594433638
$data[29,5] = 0.2487587594433639
$data[29,6] = 3
$data[29,7] = 1
$data[29,8] = 28.73046866666667
$data[29,9] = 86.191406
$data[29,10] = 0.2377611270219907
$data[29,11] = 0.2377611270219907
$data[29,12] = 63.39654681481488
$data[29,13] = 570.568921333334
$data[29,14] = 0.05914516300184646
$data[29,15] = 0.05914516300184647
$data[30,0] = 3
$data[30,1] = 1
$data[30,2] = 1.326136666666667
$data[30,3] = 3.97841
$data[30,4] = 0.1495008883450928
$data[30,5] = 0.1495008883450928
$data[30,6] = 3
$data[30,7] = 1
$data[30,8] = 3.264166666666667
$data[30,9] = 9.7925
$data[30,10] = 0.02701285365228691
$data[30,11] = 0.02701285365228691
$data[30,12] = 4.328731102777778
$data[30,13] = 38.958579925
$data[30,14] = 0.004038445617752877
$data[30,15] = 0.004038445617752877
$data[31,0] = 3
$data[31,1] = 1
$data[31,2] = 1.326136666666667
$data[31,3] = 3.97841
$data[31,4] = 0.1495008883450928
$data[31,5] = 0.1495008883450928
$data[31,6] = 3
$data[31,7] = 1
$data[31,8] = 37.12743
$data[31,9] = 111.38229
$data[31,10] = 0.3072508041078968
$data[31,11] = 0.3072508041078968
$data[31,12] = 49.2360462621
$data[31,13] = 443.1244163589
$data[31,14] = 0.04593426815887466
$data[31,15] = 0.04593426815887466
$data[32,0] = 3
$data[32,1] = 1
$data[32,2] = 1.326136666666667
$data[32,3] = 3.97841
$data[32,4] = 0.1495008883450928
$data[32,5] = 0.1495008883450928
$data[32,6] = 3
$data[32,7] = 1
$data[32,8] = 15.075014
$data[32,9] = 45.225042
$data[32,10] = 0.1247543978518794
$data[32,11] = 0.1247543978518794
$data[32,12] = 19.99152881591334
$data[32,13] = 179.92375934322
$data[32,14] = 0.0186508933038131
$data[32,15] = 0.0186508933038131
$data[33,0] = 3
$data[33,1] = 1
$data[33,2] = 1.326136666666667
$data[33,3] = 3.97841
$data[33,4] = 0.1495008883450928
$data[33,5] = 0.1495008883450928
$data[33,6] = 3
$data[33,7] = 1
$data[33,8] = 31.08386466666667
$data[33,9] = 93.25159400000001
$data[33,10] = 0.2572368303869774
$data[33,11] = 0.2572368303869774
$data[33,12] = 41.22145267617112
$data[33,13] = 370.9930740855401
$data[33,14] = 0.03845713465792908
$data[33,15] = 0.03845713465792908
$data[34,0] = 3
$data[34,1] = 1
$data[34,2] = 1.326136666666667
$data[34,3] = 3.97841
$data[34,4] = 0.1495008883450928
$data[34,5] = 0.1495008883450928
$data[34,6] = 3
$data[34,7] = 1
$data[34,8] = 5.556591666666667
$data[34,9] = 16.669775
$data[34,10] = 0.04598398697896871
$data[34,11] = 0.04598398697896871
$data[34,12] = 7.368799950861112
$data[34,13] = 66.31919955775001
$data[34,14] = 0.006874646903005001
$data[34,15] = 0.006874646903005001
$data[35,0] = 3
$data[35,1] = 1
$data[35,2] = 1.326136666666667
$data[35,3] = 3.97841
$data[35,4] = 0.1495008883450928
$data[35,5] = 0.1495008883450928
$data[35,6] = 3
$data[35,7] = 1
$data[35,8] = 28.73046866666667
$data[35,9] = 86.191406
$data[35,10] = 0.2377611270219907
$data[35,11] = 0.2377611270219907
$data[35,12] = 38.10052794938445
$data[35,13] = 342.90475154446
$data[35,14] = 0.03554549970371805
$data[35,15] = 0.03554549970371805
$ws.Range("E2:T37").Value2 = $data
Write-Host "done"